# Generate Report for Handoff
# Update the localization-status workbook to reflect the latest handoff run.
# The rows for 055501a9-...md, 0a7392f4-...md, 740d90c4-...md, ad53ac94-...md,
# b105cc69-...md and cd9a741a-...md (rows 7, 8, 10, 12, 13, 14) all shared the
# same "Latest HO/Handoff" timestamp text, so bumping that timestamp updates
# all of them at once, and each of those rows also now has a handoff type, so
# their Priority column becomes "ht".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview: Latest HO Xliff Generate Date
    $overview.Range("G$r").Value = "2016-08-28 16:22:46"

    # zh-cn: Latest Handoff Datetime
    $zhcn.Range("H$r").Value = "2016-08-28 16:22:42"

    # de-de: Latest Handoff Datetime (shares the same generate date as Overview)
    $dede.Range("H$r").Value = "2016-08-28 16:22:46"

    # zh-cn / de-de: Priority column set to "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
